$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.888.75"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").Value = "2.222.53"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "292.39"
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.60"
$ws.Range("E6").Value = "  +6.08%  "
$ws.Range("E7").Value = "  +0.60%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.35"
$ws.Range("E10").Value = "  +1.66%  "
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.50"
$ws.Range("E12").Value = "  +1.62%  "
$ws.Range("E13").Value = "  +1.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.42"
$ws.Range("E14").Value = "  +2.13%  "
$ws.Range("D15").Value = "2.562.85"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.11"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Value = "2.222.13"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.732"
$ws.Range("E18").Value = "  +1.93%  "
$ws.Range("D19").Value = "39.838.74"
$ws.Range("E19").Value = "  +1.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.51"
$ws.Range("E20").Value = "  +11.36%  "
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("E22").Value = "  +1.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.84"
$ws.Range("E23").Value = "  +1.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.90"
$ws.Range("E24").Value = "  +2.62%  "
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("E26").Value = "  +2.67%  "
$ws.Range("E27").Value = "  +0.76%  "
$ws.Range("E28").Value = "  +0.48%  "
$ws.Range("E29").Value = "  +1.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.28"
$ws.Range("E30").Value = "  +1.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.85"
$ws.Range("E31").Value = "  +2.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.50"
$ws.Range("E32").Value = "  +2.35%  "
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("E34").Value = "  +2.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0721"
$ws.Range("E35").Value = "  +3.04%  "
$ws.Range("E36").Value = "  +1.42%  "
$ws.Range("E37").Value = "  +6.12%  "
$ws.Range("E38").Value = "  +1.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.99"
$ws.Range("E39").Value = "  +0.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0993"
$ws.Range("E40").Value = "  +2.81%  "
$ws.Range("E41").Value = "  +2.61%  "
$ws.Range("D42").Value = "2.097.23"
$ws.Range("E42").Value = "  +9.54%  "
$ws.Range("E43").Value = "  +3.45%  "
$ws.Range("E44").Value = "  +5.54%  "
$ws.Range("E45").Value = "  +3.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.01"
$ws.Range("E46").Value = "  +8.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.79"
$ws.Range("E47").Value = "  +8.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.65"
$ws.Range("E48").Value = "  +1.46%  "
$ws.Range("D49").Value = "2.434.93"
$ws.Range("E49").Value = "  +0.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.01"
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "89.44"
$ws.Range("E51").Value = "  +1.35%  "
